# Update Value at Risk 95% (column N) figures for the sub category
# aggregated return rows across the Alternative, Bond and Equity sheets,
# following the new K-NN filling method.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Alternative")
$ws1.Range("N2").Value = -10.5
$ws1.Range("N3").Value = -4.51
$ws1.Range("N4").Value = -9.369999999999999
$ws1.Range("N5").Value = -81.81
$ws1.Range("N6").Value = -8.84
$ws1.Range("N7").Value = -7.24

$ws2 = $wb.Worksheets.Item("Bond")
$ws2.Range("N2").Value = -5.88
$ws2.Range("N3").Value = 3.68
$ws2.Range("N4").Value = -27.15
$ws2.Range("N5").Value = -16.35

$ws3 = $wb.Worksheets.Item("Equity")
$ws3.Range("N2").Value = -17.9
$ws3.Range("N3").Value = -15.32
$ws3.Range("N4").Value = -18.54
$ws3.Range("N5").Value = -8.81
$ws3.Range("N6").Value = -38.25
$ws3.Range("N7").Value = -16.35
$ws3.Range("N8").Value = -18.73
$ws3.Range("N9").Value = -23.04
$ws3.Range("N10").Value = -35.3
$ws3.Range("N11").Value = -1.52
